# Add two new reference columns ("hay" and "ggs") to the "grado" sheet.
# E = A + 8, F = A + 3, for every data row (2-19), with a header row
# formatted like the existing headers (copy of A1, borders removed).
# Finally make "grado" the active sheet/tab with F2:F19 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("grado")
$ws.Activate()

# Headers
$ws.Range("E1").Value = "hay"
$ws.Range("F1").Value = "ggs"

# Match header look & feel (bold Century Gothic, centered, wrapped) but
# without the vertical/side borders the original A1:D1 header row has.
$ws.Range("A1").Copy()
$ws.Range("E1:F1").PasteSpecial(-4122)
$ws.Range("E1:F1").Borders.LineStyle = -4142

# Formulas: column E = A + 8, column F = A + 3
$ws.Range("E2").Formula = "=A2+8"
$ws.Range("F2").Formula = "=A2+3"
$ws.Range("E3:E19").Formula = "=A3+8"
$ws.Range("F3:F19").Formula = "=A3+3"

# Leave the sheet selection/active tab matching the edited area.
$ws.Range("F2:F19").Select()
